# "Generate Report for Archive"
# - Update the localization status shown in the report from
#   "Ready for handoff" to "In Translation" on every sheet that reports it.
# - Re-generate (shrink) the "Latest Handoff/Target/Status" column widths on
#   each sheet to their newly auto-fit size.

$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" -------------
# Overview sheet keeps per-locale status in columns E (zh-cn) and F (de-de).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# Each locale sheet keeps its own status in column C.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- 2. Column width refresh ----------------------------------------------
# Overview: columns E and F (the two locale status columns) narrow down.
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# Locale sheets: column C (the status column) narrows down to match.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
